# Order_Assign_Test.xlsx -- "Add last source code"
#
# The diff replaces the (text-typed) values in Sheet1!AV2:AX2 with new
# ones, which in turn appends the corresponding new strings to the
# shared-strings table:
#   AV2 (FuncLoc)     : "ABCD911024" -> "ABCD354634"
#   AW2 (SAID)        : "1296412519" -> "7662871849"
#   AX2 (Previous Doc): "5908584307" -> "6119686553"
#
# All three cells hold their values as text (shared strings), even
# though AW2's cell format is numeric ("0"). Flip the format to text
# just long enough to force a text write for the purely-numeric value,
# then restore the original numeric format so the cell style is left
# exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AV2").Value = "ABCD354634"

$savedFormat = $ws.Range("AW2").NumberFormat
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value = "7662871849"
$ws.Range("AW2").NumberFormat = $savedFormat

$ws.Range("AX2").Value = "6119686553"
